$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.87"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("D3").Value = "'35.81"
$ws.Range("E3").Value = "'1.38%"
$ws.Range("D4").Value = "'5.084"
$ws.Range("E4").Value = "'0.19%"
$ws.Range("D5").Value = "'0.08082"
$ws.Range("E5").Value = "'0.88%"
$ws.Range("D6").Value = "'1.933"
$ws.Range("E6").Value = "'-0.42%"
$ws.Range("D7").Value = "'4.192"
$ws.Range("E7").Value = "'3.17%"
$ws.Range("D8").Value = "'7.761"
$ws.Range("E8").Value = "'-0.09%"
$ws.Range("D9").Value = "'0.9264"
$ws.Range("E9").Value = "'0.20%"
$ws.Range("D10").Value = "'0.1366"
$ws.Range("E10").Value = "'10.75%"
$ws.Range("D11").Value = "'0.1897"
$ws.Range("E11").Value = "'1.88%"
$ws.Range("D12").Value = "'0.09226"
$ws.Range("E12").Value = "'-4.85%"
$ws.Range("D13").Value = "'0.03432"
$ws.Range("E13").Value = "'-5.31%"
$ws.Range("D14").Value = "'0.09835"
$ws.Range("E14").Value = "'-0.15%"
$ws.Range("D15").Value = "'0.001420"
$ws.Range("E15").Value = "'2.05%"
$ws.Range("D16").Value = "'0.005771"
$ws.Range("E16").Value = "'1.16%"
$ws.Range("D17").Value = "'3.617"
$ws.Range("E17").Value = "'3.11%"
$ws.Range("D19").Value = "'0.3458"
$ws.Range("E19").Value = "'1.50%"
$ws.Range("D20").Value = "'0.1333"
$ws.Range("E20").Value = "'1.73%"
$ws.Range("D21").Value = "'4.920"
$ws.Range("E22").Value = "'-0.93%"
$ws.Range("D23").Value = "'0.04421"
$ws.Range("E23").Value = "'-2.33%"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'0.47%"
$ws.Range("D25").Value = "'0.004814"
$ws.Range("E25").Value = "'-0.40%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'3.89%"
$ws.Range("D27").Value = "'0.0003136"
$ws.Range("E27").Value = "'4.37%"
$ws.Range("D39").Value = "'0.02008"
$ws.Range("E39").Value = "'4.13%"
$ws.Range("D40").Value = "'0.04926"
$ws.Range("E40").Value = "'4.39%"
$ws.Range("D41").Value = "'0.007651"
$ws.Range("E41").Value = "'1.50%"
$ws.Range("D42").Value = "'0.01031"
$ws.Range("E42").Value = "'6.81%"
$ws.Range("D43").Value = "'0.1374"
$ws.Range("E43").Value = "'3.38%"
$ws.Range("D44").Value = "'0.002099"
$ws.Range("E44").Value = "'-0.58%"
$ws.Range("D45").Value = "'0.01103"
$ws.Range("E45").Value = "'9.57%"
$ws.Range("D46").Value = "'0.00006427"
$ws.Range("E46").Value = "'2.51%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D49").Value = "'0.001193"
$ws.Range("E49").Value = "'-19.89%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.08%"
